$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: fiscal period labels shift by one year (drop 1396/12, add 1401/12) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates shift (drop 1400-09-28, add 1402-03-11 (9)/(2)) ---
$ws.Range("D9").Value = "1400-09-29 (2)"
$ws.Range("E9").Value = "1400-10-25 (4)"
$ws.Range("F9").Value = "1401-03-07 (8)"
$ws.Range("G9").Value = "1402-03-11 (9)"
$ws.Range("H9").Value = "1402-03-11 (2)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 10897385
$ws.Range("E11").Value = 16898903
$ws.Range("F11").Value = 34281741
$ws.Range("G11").Value = 76583272
$ws.Range("H11").Value = 67160529

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---
$ws.Range("D12").Value = -7376667
$ws.Range("E12").Value = -13217035
$ws.Range("F12").Value = -20868012
$ws.Range("G12").Value = -53742202
$ws.Range("H12").Value = -45627842

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 3520718
$ws.Range("E13").Value = 3681868
$ws.Range("F13").Value = 13413729
$ws.Range("G13").Value = 22841070
$ws.Range("H13").Value = 21532687

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -79284
$ws.Range("E14").Value = -350853
$ws.Range("F14").Value = -478757
$ws.Range("G14").Value = -228237
$ws.Range("H14").Value = -544681

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) -- was "-" text, now numeric 0 ---
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 121137
$ws.Range("E16").Value = 64269
$ws.Range("F16").Value = 182133
$ws.Range("G16").Value = -240032
$ws.Range("H16").Value = 1174146

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 3562571
$ws.Range("E17").Value = 3395284
$ws.Range("F17").Value = 13117105
$ws.Range("G17").Value = 22372801
$ws.Range("H17").Value = 22162152

# --- Row 18: هزینه های مالی (Financial costs) -- D18 was "-" text, now numeric ---
$ws.Range("D18").Value = -1216
$ws.Range("E18").Value = -1345
$ws.Range("F18").Value = -6336
$ws.Range("G18").Value = -63824
$ws.Range("H18").Value = -290496

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 592483
$ws.Range("E19").Value = 491676
$ws.Range("F19").Value = 1504363
$ws.Range("G19").Value = 3217863
$ws.Range("H19").Value = 3434028

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 4153838
$ws.Range("E20").Value = 3885615
$ws.Range("F20").Value = 14615132
$ws.Range("G20").Value = 25526840
$ws.Range("H20").Value = 25305684

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -8379
$ws.Range("E21").Value = -24629
$ws.Range("F21").Value = -10437
$ws.Range("G21").Value = -30432
$ws.Range("H21").Value = -77032

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 4145459
$ws.Range("E22").Value = 3860986
$ws.Range("F22").Value = 14604695
$ws.Range("G22").Value = 25496408
$ws.Range("H22").Value = 25228652

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی -- was "-" text, now numeric 0 ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 4145459
$ws.Range("E24").Value = 3860986
$ws.Range("F24").Value = 14604695
$ws.Range("G24").Value = 25496408
$ws.Range("H24").Value = 25228652

# --- Row 25: سود هر سهم پس از کسر مالیات ---
$ws.Range("D25").Value = 921
$ws.Range("E25").Value = 858
$ws.Range("F25").Value = 1217
$ws.Range("G25").Value = 2125
$ws.Range("H25").Value = 2102

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 4500000
$ws.Range("E26").Value = 4500000
$ws.Range("F26").Value = 12000000
$ws.Range("G26").Value = 12000000
$ws.Range("H26").Value = 12000000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 118
$ws.Range("E27").Value = 110
$ws.Range("F27").Value = 417
$ws.Range("G27").Value = 728
$ws.Range("H27").Value = 721
